$d = $word.ActiveDocument

# 1) "Вариант 3.3" -> "Вариант 2.3" (only the first "3" after "Вариант " changes)
$found = $d.Content.Find.Execute("Вариант 3.3", $true, $false, $false, $false, $false, `
                                  $true, 1, $false, "Вариант 2.3", 2)

# 2) Merge the split "Красоцкий" / " М. Д." runs (with proofErr wrappers) into one run
$found2 = $d.Content.Find.Execute("Красоцкий М. Д.", $true, $false, $false, $false, $false, `
                                   $true, 1, $false, "Красоцкий М. Д.", 2)

# 3) Add <w:noProof/> to the run containing the first inline drawing
$d.InlineShapes.Item(1).Range.Font.Name = $d.InlineShapes.Item(1).Range.Font.Name
